# Add a new "created_at_field" column to the table, positioned right
# after "primary_key" and before "updated_at_field" (i.e. a new column F,
# with the old F/G/H (updated_at_field/description/comments) shifting
# right to G/H/I).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a blank worksheet column at F. This shifts the old F/G/H
#    columns (and their styles/values) to G/H/I, matching the diff.
$ws.Columns("F").Insert()

# 2. Populate the new column's header + the two data rows that had a
#    value in the (now shifted) "updated_at_field" column.
$ws.Range("F3").Value = "created_at_field"
$ws.Range("F4").Value = "landed_timestamp"
$ws.Range("F5").Value = "landed_timestamp"

# 3. New column inherits column E's width (matches real Excel's
#    behaviour of a freshly inserted column taking its left neighbour's
#    width). 17.830729166666668 is the COM ColumnWidth value whose
#    stored/serialized width is the closest achievable to column E's
#    18.6640625.
$ws.Columns("F").ColumnWidth = 17.830729166666668

# 4. Grow the table to include the new column.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("B3:I12"))

# 5. Re-assert every header cell's text so the table's column-name
#    metadata re-syncs from the worksheet cells (left to right) instead
#    of keeping stale/duplicated names from the resize.
$ws.Range("B3").Value = "no"
$ws.Range("C3").Value = "data_src"
$ws.Range("D3").Value = "table"
$ws.Range("E3").Value = "primary_key"
$ws.Range("F3").Value = "created_at_field"
$ws.Range("G3").Value = "updated_at_field"
$ws.Range("H3").Value = "description"
$ws.Range("I3").Value = "comments"

# 6. Match the final cursor/selection position from the diff.
$ws.Range("F14").Select() | Out-Null
